# Scheduled-runner update: refresh cached Universalis price snapshots
# (currentAveragePrice / NQ / HQ) and the leve profit columns that derive
# from them, per crafting-job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 5920.846
$ws.Range("I94").Value = 5920.846
$ws.Range("K94").Value = 5920.846
$ws.Range("M94").Value = -5469.846

$ws.Range("H129").Value = 962.4138
$ws.Range("I129").Value = 478.36365
$ws.Range("K129").Value = 1435.09095
$ws.Range("M129").Value = 3564.90905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8933.063
$ws.Range("I32").Value = 9146.169
$ws.Range("J32").Value = 5139.8
$ws.Range("K32").Value = 9146.169
$ws.Range("L32").Value = 5139.8
$ws.Range("M32").Value = -8859.169
$ws.Range("N32").Value = -5713.8

$ws.Range("H122").Value = 1906
$ws.Range("I122").Value = 1830.862
$ws.Range("K122").Value = 5492.586
$ws.Range("M122").Value = -3042.586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2354.5557
$ws.Range("I99").Value = 636.6667
$ws.Range("J99").Value = 3213.5
$ws.Range("K99").Value = 636.6667
$ws.Range("L99").Value = 3213.5
$ws.Range("M99").Value = 861.3333
$ws.Range("N99").Value = -6209.5

$ws.Range("H105").Value = 2451.8823
$ws.Range("I105").Value = 2455.125
$ws.Range("K105").Value = 2455.125
$ws.Range("M105").Value = -708.125

$ws.Range("H107").Value = 59183.445
$ws.Range("I107").Value = 59183.445
$ws.Range("K107").Value = 59183.445
$ws.Range("M107").Value = -57263.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2519.5
$ws.Range("I99").Value = 2648.3333
$ws.Range("J99").Value = 2133
$ws.Range("K99").Value = 2648.3333
$ws.Range("L99").Value = 2133
$ws.Range("M99").Value = -1150.3333
$ws.Range("N99").Value = -5129

$ws.Range("H126").Value = 2519.5
$ws.Range("I126").Value = 2648.3333
$ws.Range("J126").Value = 2133
$ws.Range("K126").Value = 7944.999899999999
$ws.Range("L126").Value = 6399
$ws.Range("M126").Value = -5474.999899999999
$ws.Range("N126").Value = -11339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2250
$ws.Range("J54").Value = 2250
$ws.Range("L54").Value = 6750
$ws.Range("N54").Value = -7868

$ws.Range("H117").Value = 43202.582
$ws.Range("J117").Value = 47100.09
$ws.Range("L117").Value = 141300.27
$ws.Range("N117").Value = -148184.27

$ws.Range("H118").Value = 1848.3846
$ws.Range("J118").Value = 2716.6667
$ws.Range("L118").Value = 8150.000100000001
$ws.Range("N118").Value = -10636.0001

$ws.Range("H131").Value = 850.1799999999999
$ws.Range("J131").Value = 883.98914
$ws.Range("L131").Value = 2651.96742
$ws.Range("N131").Value = -12731.96742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 24010360
$ws.Range("I11").Value = 40000000
$ws.Range("J11").Value = 25900
$ws.Range("K11").Value = 40000000
$ws.Range("L11").Value = 25900
$ws.Range("M11").Value = -39999861
$ws.Range("N11").Value = -26178

$ws.Range("H12").Value = 4151714.2
$ws.Range("I12").Value = 9334000
$ws.Range("J12").Value = 265000
$ws.Range("K12").Value = 9334000
$ws.Range("L12").Value = 265000
$ws.Range("M12").Value = -9333860
$ws.Range("N12").Value = -265280

$ws.Range("H97").Value = 18323.3
$ws.Range("I97").Value = 21727.16
$ws.Range("J97").Value = 1304
$ws.Range("K97").Value = 21727.16
$ws.Range("L97").Value = 1304
$ws.Range("M97").Value = -21231.16
$ws.Range("N97").Value = -2296

$ws.Range("H102").Value = 2185.516
$ws.Range("I102").Value = 2182.4783
$ws.Range("J102").Value = 2194.25
$ws.Range("K102").Value = 2182.4783
$ws.Range("L102").Value = 2194.25
$ws.Range("M102").Value = -560.4783000000002
$ws.Range("N102").Value = -5438.25

$ws.Range("H122").Value = 2488.625
$ws.Range("I122").Value = 2481.4
$ws.Range("J122").Value = 2524.75
$ws.Range("K122").Value = 7444.200000000001
$ws.Range("L122").Value = 7574.25
$ws.Range("M122").Value = -4994.200000000001
$ws.Range("N122").Value = -12474.25

$ws.Range("H126").Value = 2911.4285
$ws.Range("I126").Value = 2578.182
$ws.Range("J126").Value = 4133.3335
$ws.Range("K126").Value = 7734.545999999999
$ws.Range("L126").Value = 12400.0005
$ws.Range("M126").Value = -5264.545999999999
$ws.Range("N126").Value = -17340.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 28592.309
$ws.Range("J20").Value = 24700
$ws.Range("L20").Value = 24700
$ws.Range("N20").Value = -25152

$ws.Range("H93").Value = 1250.75
$ws.Range("I93").Value = 1250.75
$ws.Range("K93").Value = 1250.75
$ws.Range("M93").Value = -2.75

$ws.Range("H100").Value = 3755.1365
$ws.Range("I100").Value = 4993.5713
$ws.Range("J100").Value = 1587.875
$ws.Range("K100").Value = 4993.5713
$ws.Range("L100").Value = 1587.875
$ws.Range("M100").Value = -4452.5713
$ws.Range("N100").Value = -2669.875

$ws.Range("H132").Value = 6308
$ws.Range("I132").Value = 6894.8887
$ws.Range("J132").Value = 4547.3335
$ws.Range("K132").Value = 20684.6661
$ws.Range("L132").Value = 13642.0005
$ws.Range("M132").Value = -18154.6661
$ws.Range("N132").Value = -18702.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 52692.1
$ws.Range("I81").Value = 85787.336
$ws.Range("J81").Value = 3049.25
$ws.Range("K81").Value = 171574.672
$ws.Range("L81").Value = 6098.5
$ws.Range("M81").Value = -170513.672
$ws.Range("N81").Value = -8220.5

$ws.Range("H84").Value = 52692.1
$ws.Range("I84").Value = 85787.336
$ws.Range("J84").Value = 3049.25
$ws.Range("K84").Value = 857873.36
$ws.Range("L84").Value = 30492.5
$ws.Range("M84").Value = -852569.36
$ws.Range("N84").Value = -41100.5

$ws.Range("H122").Value = 35719790
$ws.Range("I122").Value = 50002600
$ws.Range("J122").Value = 12777.5
$ws.Range("K122").Value = 150007800
$ws.Range("L122").Value = 38332.5
$ws.Range("M122").Value = -150005350
$ws.Range("N122").Value = -43232.5

$ws.Range("H132").Value = 2864.7917
$ws.Range("I132").Value = 2149.7144
$ws.Range("J132").Value = 3865.9
$ws.Range("K132").Value = 6449.1432
$ws.Range("L132").Value = 11597.7
$ws.Range("M132").Value = -3919.1432
$ws.Range("N132").Value = -16657.7

$ws.Range("H136").Value = 1696.289
$ws.Range("I136").Value = 1769.625
$ws.Range("J136").Value = 1515.7693
$ws.Range("K136").Value = 5308.875
$ws.Range("L136").Value = 4547.3079
$ws.Range("M136").Value = -2758.875
$ws.Range("N136").Value = -9647.3079

